# Apply the "aggiornati i file excel con la durata totale dei singoli test
# di download e upload" edit:
#  - rename the two bandwidth headers to include units (Mb/s)
#  - add two new columns with the total duration (s) of the download/upload
#    tests
#  - update the chart title (profile label + file-size line) and the two
#    series display names to match the renamed headers
#  - resize/reposition the chart and widen the data columns

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- sheet data -----------------------------------------------------------

$ws.Range("B1").Value = "Banda in download (Mb/s)"
$ws.Range("C1").Value = "Banda in upload (Mb/s)"
$ws.Range("D1").Value = "Tempo totale download (s)"
$ws.Range("E1").Value = "Tempo totale upload (s)"

$ws.Range("D2").Value = 16.002
$ws.Range("E2").Value = 16.044
$ws.Range("D3").Value = 16.002
$ws.Range("E3").Value = 19.029
$ws.Range("D4").Value = 19.002
$ws.Range("E4").Value = 19.023

# --- column widths ----------------------------------------------------------
# (ColumnWidth is expressed in "characters"; the values below are chosen so
# the saved width attribute lands as close as possible to the target widths
# of 22.5 / 23.1640625 / 23.5 / 23.33203125 / 22.6640625)

$ws.Columns.Item(1).ColumnWidth = 21.666666666666668
$ws.Columns.Item(2).ColumnWidth = 22.333333333333332
$ws.Columns.Item(3).ColumnWidth = 22.666666666666668
$ws.Columns.Item(4).ColumnWidth = 22.5
$ws.Columns.Item(5).ColumnWidth = 21.833333333333332

# --- selection --------------------------------------------------------------

$ws.Range("I4").Select()

# --- chart --------------------------------------------------------------

$co = $ws.ChartObjects().Item(1)
$chart = $co.Chart

# Title: "Profilo: 30/3" + "Dimensione file: 5MB" on its own line
$chart.ChartTitle.Text = "Profilo: 30/3" + [char]10 + "Dimensione file: 5MB"

# Series display names follow the renamed headers
$chart.SeriesCollection(1).Name = "Banda in download (Mb/s)"
$chart.SeriesCollection(2).Name = "Banda in upload (Mb/s)"

# Reposition/resize the chart: anchored from A8 (col0/row7 + a bit) to
# col11/row29 (matches the target twoCellAnchor once the column widths
# above have been applied)
$co.Left = 0
$co.Top = 127
$co.Width = 1022
$co.Height = 341
